# Update "想去人数" (attendance interest count) values in the "展览" (Exhibitions)
# and "全部类型" (All Types) sheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new value for column F
$exhibitionUpdates = @{
    2  = 7070
    4  = 463
    6  = 554
    7  = 156
    8  = 121
    11 = 52
    12 = 202
    13 = 448
    15 = 1831
    17 = 3665
    21 = 25
    23 = 2294
    25 = 262
    27 = 37
    31 = 159
    32 = 1310
    33 = 113
}

$allTypesUpdates = @{
    2  = 7070
    4  = 463
    7  = 554
    8  = 156
    9  = 121
    12 = 52
    13 = 202
    14 = 448
    16 = 1831
    18 = 3665
    22 = 25
    24 = 2294
    26 = 262
    28 = 37
    32 = 159
    33 = 1310
    34 = 113
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
